$d = $word.ActiveDocument

# 1. Update activation date
$d.Content.Find.Execute(
    "Ativação: 01/01/2020", $false, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2025", 2) | Out-Null

# 2. Objetivos paragraph: add missing space and drop the trailing sentence
$d.Content.Find.Execute(
    "engenharia.Processos de metalurgia do pó metálico, de seus principais aspectos metalúrgicos, propriedades, aplicações, vantagens e desvantagens técnicas e econômicas. Identificação dos problemas comuns em componentes metálicos fundidos, soldados e sinterizados. Introdução à Manufatura Aditiva: Potencialidade e Técnicas.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "engenharia. Processos de metalurgia do pó metálico, de seus principais aspectos metalúrgicos, propriedades, aplicações, vantagens e desvantagens técnicas e econômicas. Identificação dos problemas comuns em componentes metálicos fundidos, soldados e sinterizados.",
    2) | Out-Null

# 3. Insert new docente run before Hugo Ricardo's run
$hugoRange = $d.Content
$hugoRange.Find.Execute(
    "984972 - Hugo Ricardo Zschommler Sandim", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$insertRange = $hugoRange.Duplicate
$insertRange.Collapse(1)
$insertRange.InsertBefore("3586455 - Cassius Olivio Figueiredo Terra Ruchert`v")

# 4. Programa resumido bullet paragraph
$d.Content.Find.Execute(
    "1. TÉCNICAS DE JUNÇÃO DE MATERIAIS; 2. PROCESSOS DE SOLDAGEM; 3. NOMENCLATURA DAS JUNTAS SOLDADAS; 4. METALURGIA FÍSICA DAS REGIÕES SOLDADAS; 5. SEGURANÇA NO PROCESSO DE SOLDAGEM; 6. APLICAÇÕES DE JUNTAS SOLDADAS EM ENGENHARIA; 7.PÓS METALICOS – OBTENÇÃO, CARACTERIZAÇÃO E APLICAÇÃO NA METALURGIA DO PÓ. 8. TÉCNICAS DE MISTURA, 9. PROCESSOS DE FABRICAÇÃO DE PEÇAS VERDES, 10. SINTERIZAÇÃO, 11. UTILIZAÇÃO DO LASER E DE FEIXE DE ELÉTRONS12. PRÁTICA EXPERIMENTAL SUPERVISIONADA.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "1. TÉCNICAS DE JUNÇÃO DE MATERIAIS; 2. PROCESSOS DE SOLDAGEM; 3. NOMENCLATURA DAS JUNTAS SOLDADAS; 4. METALURGIA FÍSICA DAS REGIÕES SOLDADAS; 5. SEGURANÇA NO PROCESSO DE SOLDAGEM; 6. APLICAÇÕES DE JUNTAS SOLDADAS EM ENGENHARIA; 7.PÓS METALICOS – OBTENÇÃO, CARACTERIZAÇÃO E APLICAÇÃO NA METALURGIA DO PÓ. 8. TÉCNICAS DE MISTURA, 9. PROCESSOS DE FABRICAÇÃO DE PEÇAS VERDES, 10. SINTERIZAÇÃO, 11. MANUFATURA ADITIVA, 12. UTILIZAÇÃO DE LASER E FEIXE DE ELÉTRONS, 13. PARÂMETROS RELEVANTES, 14. ESTUDOS DE CASOS APLICADOS.",
    2) | Out-Null

# 5. Full "Programa" paragraph (use Range.Text= so literal straight quotes
# survive verbatim — Find.Execute's replacement text goes through
# smart-quote autocorrection, which would corrupt the "verdes" quotes).
$progRange = $d.Content
$progRange.Find.Execute(
    "1. TÉCNICAS DE JUNÇÃO DE MATERIAIS: Razões técnicas para a junção de materiais, junção por difusão, brasagem, soldagem por explosão, elementos de fixação. 2. PROCESSOS DE SOLDAGEM: Definição de soldagem por fusão, física da soldagem, principais processos de soldagem. 3. NOMENCLATURA DAS JUNTAS SOLDADAS: Desenho e simbologia para soldagem, símbolos básicos, tipos de juntas e soldas, simbologia para soldas em desenho. 4. METALURGIA FÍSICA DAS REGIÕES SOLDADAS: metalurgia da soldagem, estruturas de solidificação, transformações de fase pós-soldagem, transformações de fases em juntas de aço soldadas, ligas de alumínio, ligas de cobre e em metais e ligas especiais. 5. SEGURANÇA NO PROCESSO DE SOLDAGEM: Problemas associados à vaporização de metais, luminosidade, calor e eletricidade. 6. APLICAÇÕES DE JUNTAS SOLDADAS EM ENGENHARIA: Exemplos de estruturas soldadas em engenharia, descontinuidades e defeitos de soldagem, métodos de inspeção em soldas. 7. PRÁTICA EXPERIMENTAL SUPERVISIONADA: Caracterização microestrutural de juntas soldadas (materiais e processos a serem definidos na ocasião da prática experimental), incluindo a redação de relatório técnico de cada grupo. 8. Pós Metálicos - obtenção por processos químicos, termoquímicos, eletrolíticos, atomização e moagem, Caracterização de pós e sua aplicação na metalurgia do pó.9. Técnicas de mistura, aspectos sobre o transporte e armazenamento de pós, 10-Processos de fabricação de peças `"verdes`" por compactação uniaxial e isostática, 11- Técnicas de sinterização e fenômenos envolvidos, 12 Sinterização/refusão a LASER para prototipagem rápida (impressão 3D). Feixe de elétrons: obtenção e aplicações.",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null

$progRange.Text = "Programa1. TÉCNICAS DE JUNÇÃO DE MATERIAIS: Razões técnicas para a junção de materiais, junção por difusão, brasagem, soldagem por explosão, elementos de fixação. 2. PROCESSOS DE SOLDAGEM: Definição de soldagem por fusão, física da soldagem, principais processos de soldagem. 3. NOMENCLATURA DAS JUNTAS SOLDADAS: Desenho e simbologia para soldagem, símbolos básicos, tipos de juntas e soldas, simbologia para soldas em desenho. 4. METALURGIA FÍSICA DAS REGIÕES SOLDADAS: metalurgia da soldagem, estruturas de solidificação, transformações de fase pós-soldagem, transformações de fases em juntas de aço soldadas, ligas de alumínio, ligas de cobre e em metais e ligas especiais. 5. SEGURANÇA NO PROCESSO DE SOLDAGEM: Problemas associados à vaporização de metais, luminosidade, calor e eletricidade. 6. APLICAÇÕES DE JUNTAS SOLDADAS EM ENGENHARIA: Exemplos de estruturas soldadas em engenharia, descontinuidades e defeitos de soldagem, métodos de inspeção em soldas. 7. PRÁTICA EXPERIMENTAL SUPERVISIONADA: Caracterização microestrutural de juntas soldadas (materiais e processos a serem definidos na ocasião da prática experimental), incluindo a redação de relatório técnico de cada grupo. 8. Pós Metálicos - obtenção por processos químicos, termoquímicos, eletrolíticos, atomização e moagem, Caracterização de pós e sua aplicação na metalurgia do pó. 9. Técnicas de mistura, aspectos sobre o transporte e armazenamento de pós, 10. Processos de fabricação de peças `"verdes`" por compactação uniaxial e isostática, 11. Técnicas de sinterização e fenômenos envolvidos, 12. Manufatura aditiva (impressão 3D). 13. Fontes de calor (laser e feixe de elétrons: obtenção e aplicações), características desejáveis dos pós, parâmetros relevantes no processamento a laser e por feixe de elétrons. 14. Estudos de casos e comparação entre as técnicas estudadas no semestre."

# 6. Método
$d.Content.Find.Execute(
    "O aluno será avaliado por duas avaliações, sendo que a segunda avaliação terá peso 2.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Duas avaliações escritas, compostas por provas e que poderão ser complementadas por trabalhos ou relatórios de experimentos realizados em laboratório.",
    2) | Out-Null

# 7. Critério
$d.Content.Find.Execute(
    "Nota Final NF = [Avaliação 1 + 2*(Avaliação 2)]/3",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "A cada avaliação (compreendendo uma prova, complementada por trabalho ou relatório) será atribuído grau entre zero e dez.",
    2) | Out-Null

# 8. Norma de recuperação
$d.Content.Find.Execute(
    "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Avaliação escrita. Para aprovação, a média entre a avaliação de Recuperação e o grau obtido no semestre deve ser maior ou igual a cinco.",
    2) | Out-Null

Write-Output "done"
